{"js": "// Each entry is [oldText, newText] for the date line plus every\n// '<addend>+<addend>=' / '<minuend>-<subtrahend>=' table cell, in document order.\nconst pairs = [\n  [\"2024-05-22 Wednesday\", \"2024-05-23 Thursday\"],\n  [\"19+39=\", \"11+32=\"],\n  [\"53+14=\", \"47-28=\"],\n  [\"48-3=\", \"4+50=\"],\n  [\"31+27=\", \"98-91=\"],\n  [\"39+21=\", \"49+7=\"],\n  [\"49+27=\", \"96-4=\"],\n  [\"66-38=\", \"0+91=\"],\n  [\"39-20=\", \"7+56=\"],\n  [\"14+17=\", \"20+59=\"],\n  [\"11+38=\", \"64-16=\"],\n  [\"51-14=\", \"44+30=\"],\n  [\"76-17=\", \"55-18=\"],\n  [\"28+70=\", \"2+81=\"],\n  [\"55-9=\", \"13+48=\"],\n  [\"40+12=\", \"65-56=\"],\n  [\"81+7=\", \"90-7=\"],\n  [\"63+5=\", \"16+19=\"],\n  [\"96-41=\", \"24-14=\"],\n  [\"43+53=\", \"34-9=\"],\n  [\"47-12=\", \"75-61=\"],\n  [\"9+33=\", \"81-0=\"],\n  [\"39+12=\", \"49+30=\"],\n  [\"40+8=\", \"50+1=\"],\n  [\"86-15=\", \"2+9=\"],\n  [\"52+45=\", \"8+32=\"],\n  [\"99-68=\", \"83+13=\"],\n  [\"20+73=\", \"18+34=\"],\n  [\"22+20=\", \"82-56=\"],\n  [\"74-72=\", \"72-54=\"],\n  [\"81-46=\", \"60-49=\"],\n  [\"32+67=\", \"24+52=\"],\n  [\"16+21=\", \"72-33=\"],\n  [\"42-39=\", \"85-59=\"],\n  [\"73+8=\", \"96-52=\"],\n  [\"23+30=\", \"23-9=\"],\n  [\"25+34=\", \"10+67=\"],\n  [\"71-57=\", \"73-23=\"],\n  [\"53-39=\", \"37-19=\"],\n  [\"76-60=\", \"74+24=\"],\n  [\"9+51=\", \"55+24=\"],\n  [\"37+50=\", \"3+48=\"],\n  [\"20-4=\", \"19+37=\"],\n  [\"41-41=\", \"80-3=\"],\n  [\"25+70=\", \"94-9=\"],\n  [\"33-10=\", \"51-37=\"],\n  [\"47+8=\", \"91-61=\"],\n  [\"97-79=\", \"7+16=\"],\n  [\"64-14=\", \"25+41=\"],\n  [\"14+37=\", \"18+6=\"],\n  [\"22-12=\", \"31+17=\"],\n  [\"54-35=\", \"24+38=\"],\n  [\"68+15=\", \"59-46=\"],\n  [\"75-52=\", \"39-19=\"],\n  [\"57+6=\", \"1+74=\"],\n  [\"74-2=\", \"62-36=\"],\n  [\"84+14=\", \"57+1=\"],\n  [\"91-71=\", \"30+48=\"],\n  [\"63+19=\", \"71-6=\"],\n  [\"58-51=\", \"0+15=\"],\n  [\"73-30=\", \"14+68=\"],\n  [\"73-46=\", \"88-4=\"],\n  [\"41+42=\", \"92-32=\"],\n  [\"55-6=\", \"53+28=\"],\n  [\"27+58=\", \"93-10=\"],\n  [\"82-24=\", \"73+19=\"],\n  [\"31+24=\", \"6-2=\"],\n  [\"22+5=\", \"33+19=\"],\n  [\"24-20=\", \"86-4=\"],\n  [\"57-41=\", \"88-35=\"],\n  [\"97-90=\", \"94-3=\"],\n  [\"71+4=\", \"36+6=\"],\n  [\"9+39=\", \"62-18=\"],\n  [\"40+25=\", \"26+7=\"],\n  [\"32+44=\", \"29+14=\"],\n  [\"18+76=\", \"26+16=\"],\n  [\"91-73=\", \"4+53=\"],\n  [\"21+19=\", \"20+24=\"],\n  [\"34+14=\", \"32+49=\"],\n  [\"87-51=\", \"79-0=\"],\n  [\"29+45=\", \"84+6=\"],\n  [\"80-2=\", \"45+4=\"],\n  [\"44+0=\", \"13+11=\"],\n  [\"20+39=\", \"18+4=\"],\n  [\"29-3=\", \"92-7=\"],\n  [\"68-5=\", \"48-7=\"],\n  [\"24-13=\", \"36-31=\"],\n  [\"88-12=\", \"90-36=\"],\n  [\"87-37=\", \"27+21=\"],\n  [\"38+20=\", \"11+5=\"],\n  [\"43-34=\", \"91-6=\"],\n  [\"65-37=\", \"6+59=\"],\n  [\"77-23=\", \"34+45=\"],\n  [\"98-20=\", \"76-70=\"],\n  [\"67-66=\", \"85-6=\"],\n  [\"4+8=\", \"4+79=\"],\n  [\"61-60=\", \"80+17=\"],\n  [\"61-37=\", \"17+18=\"],\n  [\"62-31=\", \"53+35=\"],\n  [\"38+11=\", \"33+52=\"],\n  [\"15+7=\", \"71+18=\"]\n];\n\n// Replace each occurrence in place: find the unique run containing the old\n// text and overwrite it with the new text, preserving its formatting.\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Each pair is (oldText, newText) for the date line plus every\n# '<addend>+<addend>=' / '<minuend>-<subtrahend>=' table cell, in document order.\n$d = $word.ActiveDocument\n$pairs = @(\n  @('2024-05-22 Wednesday', '2024-05-23 Thursday'),\n  @('19+39=', '11+32='),\n  @('53+14=', '47-28='),\n  @('48-3=', '4+50='),\n  @('31+27=', '98-91='),\n  @('39+21=', '49+7='),\n  @('49+27=', '96-4='),\n  @('66-38=', '0+91='),\n  @('39-20=', '7+56='),\n  @('14+17=', '20+59='),\n  @('11+38=', '64-16='),\n  @('51-14=', '44+30='),\n  @('76-17=', '55-18='),\n  @('28+70=', '2+81='),\n  @('55-9=', '13+48='),\n  @('40+12=', '65-56='),\n  @('81+7=', '90-7='),\n  @('63+5=', '16+19='),\n  @('96-41=', '24-14='),\n  @('43+53=', '34-9='),\n  @('47-12=', '75-61='),\n  @('9+33=', '81-0='),\n  @('39+12=', '49+30='),\n  @('40+8=', '50+1='),\n  @('86-15=', '2+9='),\n  @('52+45=', '8+32='),\n  @('99-68=', '83+13='),\n  @('20+73=', '18+34='),\n  @('22+20=', '82-56='),\n  @('74-72=', '72-54='),\n  @('81-46=', '60-49='),\n  @('32+67=', '24+52='),\n  @('16+21=', '72-33='),\n  @('42-39=', '85-59='),\n  @('73+8=', '96-52='),\n  @('23+30=', '23-9='),\n  @('25+34=', '10+67='),\n  @('71-57=', '73-23='),\n  @('53-39=', '37-19='),\n  @('76-60=', '74+24='),\n  @('9+51=', '55+24='),\n  @('37+50=', '3+48='),\n  @('20-4=', '19+37='),\n  @('41-41=', '80-3='),\n  @('25+70=', '94-9='),\n  @('33-10=', '51-37='),\n  @('47+8=', '91-61='),\n  @('97-79=', '7+16='),\n  @('64-14=', '25+41='),\n  @('14+37=', '18+6='),\n  @('22-12=', '31+17='),\n  @('54-35=', '24+38='),\n  @('68+15=', '59-46='),\n  @('75-52=', '39-19='),\n  @('57+6=', '1+74='),\n  @('74-2=', '62-36='),\n  @('84+14=', '57+1='),\n  @('91-71=', '30+48='),\n  @('63+19=', '71-6='),\n  @('58-51=', '0+15='),\n  @('73-30=', '14+68='),\n  @('73-46=', '88-4='),\n  @('41+42=', '92-32='),\n  @('55-6=', '53+28='),\n  @('27+58=', '93-10='),\n  @('82-24=', '73+19='),\n  @('31+24=', '6-2='),\n  @('22+5=', '33+19='),\n  @('24-20=', '86-4='),\n  @('57-41=', '88-35='),\n  @('97-90=', '94-3='),\n  @('71+4=', '36+6='),\n  @('9+39=', '62-18='),\n  @('40+25=', '26+7='),\n  @('32+44=', '29+14='),\n  @('18+76=', '26+16='),\n  @('91-73=', '4+53='),\n  @('21+19=', '20+24='),\n  @('34+14=', '32+49='),\n  @('87-51=', '79-0='),\n  @('29+45=', '84+6='),\n  @('80-2=', '45+4='),\n  @('44+0=', '13+11='),\n  @('20+39=', '18+4='),\n  @('29-3=', '92-7='),\n  @('68-5=', '48-7='),\n  @('24-13=', '36-31='),\n  @('88-12=', '90-36='),\n  @('87-37=', '27+21='),\n  @('38+20=', '11+5='),\n  @('43-34=', '91-6='),\n  @('65-37=', '6+59='),\n  @('77-23=', '34+45='),\n  @('98-20=', '76-70='),\n  @('67-66=', '85-6='),\n  @('4+8=', '4+79='),\n  @('61-60=', '80+17='),\n  @('61-37=', '17+18='),\n  @('62-31=', '53+35='),\n  @('38+11=', '33+52='),\n  @('15+7=', '71+18=')\n)\n# For every pair, find the unique occurrence of the old text in the document\n# and replace it with the new text (wdReplaceAll is safe here since each old\n# string occurs exactly once), preserving the run's formatting.\nforeach ($p in $pairs) {\n    $r = $d.Content\n    $found = $r.Find.Execute($p[0], $true, $false, $false, $false, $false, $true, 1, $false, $p[1], 2)\n    if (-not $found) {\n        Write-Output (\"NOT FOUND: \" + $p[0])\n    }\n}\nWrite-Output \"done\""}
